$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text so values like "66.32" are not
# auto-converted to numbers (and to preserve exact string formatting).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '37.500.33'
$ws.Range("E2").Value = '  +5.75%  '
$ws.Range("D3").Value = '2.058.17'
$ws.Range("E3").Value = '  +4.38%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '252.34'
$ws.Range("E5").Value = '  +3.34%  '
$ws.Range("D6").Value = '0.650'
$ws.Range("E6").Value = '  +2.53%  '
$ws.Range("D7").Value = '66.32'
$ws.Range("E7").Value = '  +16.67%  '
$ws.Range("E8").Value = '  -0.01%  '
$ws.Range("E9").Value = '  +6.77%  '
$ws.Range("D10").Value = '59.59'
$ws.Range("E10").Value = '  +1.19%  '
$ws.Range("E11").Value = '  +5.40%  '
$ws.Range("E12").Value = '  +1.78%  '
$ws.Range("E13").Value = '  -3.26%  '
$ws.Range("D14").Value = '14.91'
$ws.Range("E14").Value = '  +5.46%  '
$ws.Range("D15").Value = '2.360.29'
$ws.Range("E15").Value = '  +4.42%  '
$ws.Range("D16").Value = '21.37'
$ws.Range("E16").Value = '  +22.74%  '
$ws.Range("D17").Value = '5.59'
$ws.Range("E17").Value = '  +6.67%  '
$ws.Range("D18").Value = '2.043.71'
$ws.Range("E18").Value = '  +3.64%  '
$ws.Range("D19").Value = '37.233.51'
$ws.Range("E19").Value = '  +5.17%  '
$ws.Range("D20").Value = '73.96'
$ws.Range("E20").Value = '  +3.93%  '
$ws.Range("D21").Value = '0.0₃0879'
$ws.Range("E21").Value = '  +4.89%  '
$ws.Range("D22").Value = '5.48'
$ws.Range("E22").Value = '  +7.07%  '
$ws.Range("D23").Value = '240.07'
$ws.Range("E23").Value = '  +3.35%  '
$ws.Range("E24").Value = '  +5.89%  '
$ws.Range("E25").Value = '  -0.11%  '
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").Value = '9.72'
$ws.Range("E27").Value = '  +7.36%  '
$ws.Range("D28").Value = '161.63'
$ws.Range("E28").Value = '  -0.94%  '
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("E30").Value = '  +9.28%  '
$ws.Range("E32").Value = '  +24.56%  '
$ws.Range("D33").Value = '1.20'
$ws.Range("E33").Value = '  +6.47%  '
$ws.Range("E34").Value = '  +11.92%  '
$ws.Range("D35").Value = '0.0622'
$ws.Range("E35").Value = '  +5.34%  '
$ws.Range("D36").Value = '2.45'
$ws.Range("E36").Value = '  +4.72%  '
$ws.Range("E37").Value = '  +0.03%  '
$ws.Range("E38").Value = '  +4.69%  '
$ws.Range("E39").Value = '  +18.93%  '
$ws.Range("D40").Value = '3.07'
$ws.Range("E40").Value = '  +36.52%  '
$ws.Range("E41").Value = '  +17.57%  '
$ws.Range("E42").Value = '  +3.30%  '
$ws.Range("E43").Value = '  +5.01%  '
$ws.Range("D44").Value = '1.15'
$ws.Range("E44").Value = '  +6.76%  '
$ws.Range("D46").Value = '17.04'
$ws.Range("E46").Value = '  +7.42%  '
$ws.Range("D47").Value = '95.66'
$ws.Range("E47").Value = '  +4.98%  '
$ws.Range("D48").Value = '7.94'
$ws.Range("E48").Value = '  +6.06%  '
$ws.Range("D49").Value = '1.417.96'
$ws.Range("E49").Value = '  +2.96%  '
$ws.Range("E50").Value = '  +2.40%  '
$ws.Range("D51").Value = '46.83'
$ws.Range("E51").Value = '  +3.02%  '

# Restore the original (default) cell style now that the text values are set.
$priceRange.Style = "Normal"

